$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Col1a1"
$ws.Range("C2").Value = "Itga2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 10.379369
$ws.Range("H2").Value = 31.138107
$ws.Range("I2").Value = 0.01614698522449884
$ws.Range("J2").Value = 0.01614698522449883
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.888921
$ws.Range("N2").Value = 5.666763
$ws.Range("O2").Value = 0.4551532417350329
$ws.Range("P2").Value = 0.4551532417350328
$ws.Range("Q2").Value = 19.605808070849
$ws.Range("R2").Value = 176.452272637641
$ws.Range("S2").Value = 0.007349352669178322
$ws.Range("T2").Value = 0.007349352669178319

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Col1a1"
$ws.Range("C3").Value = "Itga2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 10.379369
$ws.Range("H3").Value = 31.138107
$ws.Range("I3").Value = 0.01614698522449884
$ws.Range("J3").Value = 0.01614698522449883
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.087098333333333
$ws.Range("N3").Value = 3.261295
$ws.Range("O3").Value = 0.2619465454094788
$ws.Range("P3").Value = 0.2619465454094787
$ws.Range("Q3").Value = 11.28339474095167
$ws.Range("R3").Value = 101.550552668565
$ws.Range("S3").Value = 0.004229646998335367
$ws.Range("T3").Value = 0.004229646998335365

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Col1a1"
$ws.Range("C4").Value = "Itga2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 10.379369
$ws.Range("H4").Value = 31.138107
$ws.Range("I4").Value = 0.01614698522449884
$ws.Range("J4").Value = 0.01614698522449883
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.174057666666667
$ws.Range("N4").Value = 3.522173
$ws.Range("O4").Value = 0.2829002128554884
$ws.Range("P4").Value = 0.2829002128554884
$ws.Range("Q4").Value = 12.18597774961233
$ws.Range("R4").Value = 109.673799746511
$ws.Range("S4").Value = 0.004567985556985146
$ws.Range("T4").Value = 0.004567985556985146

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Col1a1"
$ws.Range("C5").Value = "Itga2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 604.0312093333333
$ws.Range("H5").Value = 1812.093628
$ws.Range("I5").Value = 0.9396797639857967
$ws.Range("J5").Value = 0.9396797639857967
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 1.888921
$ws.Range("N5").Value = 5.666763
$ws.Range("O5").Value = 0.4551532417350329
$ws.Range("P5").Value = 0.4551532417350328
$ws.Range("Q5").Value = 1140.967235965129
$ws.Range("R5").Value = 10268.70512368616
$ws.Range("S5").Value = 0.427698290770946
$ws.Range("T5").Value = 0.4276982907709459

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Col1a1"
$ws.Range("C6").Value = "Itga2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 604.0312093333333
$ws.Range("H6").Value = 1812.093628
$ws.Range("I6").Value = 0.9396797639857967
$ws.Range("J6").Value = 0.9396797639857967
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.087098333333333
$ws.Range("N6").Value = 3.261295
$ws.Range("O6").Value = 0.2619465454094788
$ws.Range("P6").Value = 0.2619465454094787
$ws.Range("Q6").Value = 656.6413209475844
$ws.Range("R6").Value = 5909.771888528259
$ws.Range("S6").Value = 0.2461458679672738
$ws.Range("T6").Value = 0.2461458679672737

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Col1a1"
$ws.Range("C7").Value = "Itga2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 604.0312093333333
$ws.Range("H7").Value = 1812.093628
$ws.Range("I7").Value = 0.9396797639857967
$ws.Range("J7").Value = 0.9396797639857967
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.174057666666667
$ws.Range("N7").Value = 3.522173
$ws.Range("O7").Value = 0.2829002128554884
$ws.Range("P7").Value = 0.2829002128554884
$ws.Range("Q7").Value = 709.1674722237382
$ws.Range("R7").Value = 6382.507250013644
$ws.Range("S7").Value = 0.265835605247577
$ws.Range("T7").Value = 0.265835605247577

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Col1a1"
$ws.Range("C8").Value = "Itga2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 28.39480333333333
$ws.Range("H8").Value = 85.18441
$ws.Range("I8").Value = 0.04417325078970442
$ws.Range("J8").Value = 0.04417325078970442
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 1.888921
$ws.Range("N8").Value = 5.666763
$ws.Range("O8").Value = 0.4551532417350329
$ws.Range("P8").Value = 0.4551532417350328
$ws.Range("Q8").Value = 53.63554030720333
$ws.Range("R8").Value = 482.7198627648299
$ws.Range("S8").Value = 0.02010559829490857
$ws.Range("T8").Value = 0.02010559829490856

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Col1a1"
$ws.Range("C9").Value = "Itga2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 28.39480333333333
$ws.Range("H9").Value = 85.18441
$ws.Range("I9").Value = 0.04417325078970442
$ws.Range("J9").Value = 0.04417325078970442
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.087098333333333
$ws.Range("N9").Value = 3.261295
$ws.Range("O9").Value = 0.2619465454094788
$ws.Range("P9").Value = 0.2619465454094787
$ws.Range("Q9").Value = 30.86794337899444
$ws.Range("R9").Value = 277.8114904109499
$ws.Range("S9").Value = 0.0115710304438696
$ws.Range("T9").Value = 0.0115710304438696

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Col1a1"
$ws.Range("C10").Value = "Itga2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 28.39480333333333
$ws.Range("H10").Value = 85.18441
$ws.Range("I10").Value = 0.04417325078970442
$ws.Range("J10").Value = 0.04417325078970442
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.174057666666667
$ws.Range("N10").Value = 3.522173
$ws.Range("O10").Value = 0.2829002128554884
$ws.Range("P10").Value = 0.2829002128554884
$ws.Range("Q10").Value = 33.33713654699222
$ws.Range("R10").Value = 300.03422892293
$ws.Range("S10").Value = 0.01249662205092625
$ws.Range("T10").Value = 0.01249662205092625
